$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hours worked on 2019-07-28 (row 21) from 1 to 2
$ws.Range("H21").Value = 2

# Move the active selection to H21 (matches the saved cursor position)
$ws.Range("H21").Select()
